# Apply the "Office Theme" colour scheme to the presentation's (single
# reachable) theme -- corresponds to the clrScheme swap recorded for
# ppt/theme/theme1.xml in the target diff (Integral / Red-Violet ->
# Office Theme colours; the font/format schemes were already identical
# between the two theme parts, so only the 12 colour slots change).
$p = $ppt.ActivePresentation
$master = $p.SlideMaster
$theme = $master.Theme
$colors = $theme.ThemeColorScheme

# Index order for ThemeColorScheme matches the OOXML <a:clrScheme> child
# order: dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink. RGB values use
# the OLE BGR-packed integer (R | G<<8 | B<<16), i.e. VBA's RGB(r,g,b).
$colors.Item(1).RGB  = 0         # dk1      000000
$colors.Item(2).RGB  = 16777215  # lt1      FFFFFF
$colors.Item(3).RGB  = 6968388   # dk2      44546A
$colors.Item(4).RGB  = 15132391  # lt2      E7E6E6
$colors.Item(5).RGB  = 13998939  # accent1  5B9BD5
$colors.Item(6).RGB  = 3243501   # accent2  ED7D31
$colors.Item(7).RGB  = 10855845  # accent3  A5A5A5
$colors.Item(8).RGB  = 49407     # accent4  FFC000
$colors.Item(9).RGB  = 12874308  # accent5  4472C4
$colors.Item(10).RGB = 4697456   # accent6  70AD47
$colors.Item(11).RGB = 12673797  # hlink    0563C1
$colors.Item(12).RGB = 7491477   # folHlink 954F72

# Slide 5 ("B1- TYPES OF FINANCIAL DOCUMENTS") has a table whose style
# was switched to a different built-in table style GUID.
$slide = $p.Slides.Item(5)
for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
    $shape = $slide.Shapes.Item($i)
    if ($shape.HasTable) {
        $shape.Table.ApplyStyle("{B6E79DCC-191E-4DFC-975E-CA08B02663E9}")
    }
}
